# Apply updated cryptocurrency price/volume data (and the TRON/WrappedEther
# row swap) as of the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'63.146.98"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "'  -2.60%  "
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(3, 4).Value = "'3.238.25"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "'  -2.69%  "
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(4, 4).Value = "'1.00"
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = "'  +0.03%  "
$ws.Cells.Item(4, 5).Style = "Normal"
$ws.Cells.Item(5, 4).Value = "'517.14"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "'  -2.76%  "
$ws.Cells.Item(5, 5).Style = "Normal"
$ws.Cells.Item(6, 4).Value = "'168.44"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "'  -7.92%  "
$ws.Cells.Item(6, 5).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "'  -4.97%  "
$ws.Cells.Item(7, 5).Style = "Normal"
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "'  -0.02%  "
$ws.Cells.Item(8, 5).Style = "Normal"
$ws.Cells.Item(9, 4).Value = "'3.230.89"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "'  -2.81%  "
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "'  -5.69%  "
$ws.Cells.Item(10, 5).Style = "Normal"
$ws.Cells.Item(11, 4).Value = "'51.63"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "'  -12.67%  "
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "'  -3.26%  "
$ws.Cells.Item(12, 5).Style = "Normal"
$ws.Cells.Item(13, 4).Value = "'0.0000251"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "'  -5.21%  "
$ws.Cells.Item(13, 5).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "'  -4.58%  "
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(15, 4).Value = "'3.743.78"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "'  -2.65%  "
$ws.Cells.Item(15, 5).Style = "Normal"
$ws.Cells.Item(16, 2).Value = "'WrappedEther"
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 3).Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).Value = "'3.239.18"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "'  -2.32%  "
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(17, 2).Value = "'TRON"
$ws.Cells.Item(17, 2).Style = "Normal"
$ws.Cells.Item(17, 3).Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(17, 3).Style = "Normal"
$ws.Cells.Item(17, 4).Value = "'0.114"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "'  -2.78%  "
$ws.Cells.Item(17, 5).Style = "Normal"
$ws.Cells.Item(18, 4).Value = "'62.950.26"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "'  -2.62%  "
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(19, 4).Value = "'16.99"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "'  -4.90%  "
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(20, 4).Value = "'10.90"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "'  -2.65%  "
$ws.Cells.Item(20, 5).Style = "Normal"
$ws.Cells.Item(21, 4).Value = "'0.935"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "'  -3.23%  "
$ws.Cells.Item(21, 5).Style = "Normal"
$ws.Cells.Item(22, 4).Value = "'368.56"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "'  -2.05%  "
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(23, 4).Value = "'4.12"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "'  +4.95%  "
$ws.Cells.Item(23, 5).Style = "Normal"
$ws.Cells.Item(24, 4).Value = "'11.07"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "'  -2.35%  "
$ws.Cells.Item(24, 5).Style = "Normal"
$ws.Cells.Item(25, 4).Value = "'79.98"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "'  -1.71%  "
$ws.Cells.Item(25, 5).Style = "Normal"
$ws.Cells.Item(26, 4).Value = "'3.61"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "'  -6.23%  "
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "'  +0.34%  "
$ws.Cells.Item(27, 5).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "'  -3.17%  "
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(29, 4).Value = "'10.98"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "'  -5.84%  "
$ws.Cells.Item(29, 5).Style = "Normal"
$ws.Cells.Item(30, 4).Value = "'7.91"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "'  -6.97%  "
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(31, 4).Value = "'27.89"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "'  -4.42%  "
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(32, 4).Value = "'615.30"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "'  -4.79%  "
$ws.Cells.Item(32, 5).Style = "Normal"
$ws.Cells.Item(33, 4).Value = "'6.48"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "'  -5.45%  "
$ws.Cells.Item(33, 5).Style = "Normal"
$ws.Cells.Item(34, 4).Value = "'10.95"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "'  -3.80%  "
$ws.Cells.Item(34, 5).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "'  -3.03%  "
$ws.Cells.Item(35, 5).Style = "Normal"
$ws.Cells.Item(36, 4).Value = "'56.54"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "'  -5.93%  "
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "'  -0.03%  "
$ws.Cells.Item(37, 5).Style = "Normal"
$ws.Cells.Item(38, 4).Value = "'34.93"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "'  -5.63%  "
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(39, 4).Value = "'0.369"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "'  -7.25%  "
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(40, 4).Value = "'0.999"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "'  -0.05%  "
$ws.Cells.Item(40, 5).Style = "Normal"
$ws.Cells.Item(41, 4).Value = "'0.0₃0704"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "'  -4.10%  "
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(42, 4).Value = "'2.59"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "'  +1.72%  "
$ws.Cells.Item(42, 5).Style = "Normal"
$ws.Cells.Item(43, 4).Value = "'0.122"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "'  -4.54%  "
$ws.Cells.Item(43, 5).Style = "Normal"
$ws.Cells.Item(44, 4).Value = "'3.03"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "'  +3.40%  "
$ws.Cells.Item(44, 5).Style = "Normal"
$ws.Cells.Item(45, 4).Value = "'2.810.31"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "'  -3.60%  "
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "'  -0.55%  "
$ws.Cells.Item(46, 5).Style = "Normal"
$ws.Cells.Item(47, 4).Value = "'0.0387"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "'  -4.57%  "
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "'  -7.62%  "
$ws.Cells.Item(48, 5).Style = "Normal"
$ws.Cells.Item(49, 4).Value = "'2.96"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "'  -0.80%  "
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(50, 4).Value = "'136.19"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "'  +1.76%  "
$ws.Cells.Item(50, 5).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "'  -3.98%  "
$ws.Cells.Item(51, 5).Style = "Normal"
